$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric must be forced to Text format
# so Excel keeps them as strings (preserving formatting like trailing zeros).
$textCells = @("D5","D6","D8","D10","D16","D19","D20","D21","D22","D24","D25","D30","D32","D33","D34","D35","D36","D37","D40","D41","D44","D46","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "58.328.75"
$ws.Range("E2").Value = "  -5.01%  "
$ws.Range("D3").Value = "2.565.24"
$ws.Range("E3").Value = "  -4.43%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "507.24"
$ws.Range("E5").Value = "  -5.03%  "
$ws.Range("D6").Value = "145.13"
$ws.Range("E6").Value = "  -7.83%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").Value = "0.575"
$ws.Range("E8").Value = "  -2.83%  "
$ws.Range("D9").Value = "2.578.35"
$ws.Range("E9").Value = "  -4.80%  "
$ws.Range("D10").Value = "6.29"
$ws.Range("E10").Value = "  -4.91%  "
$ws.Range("E11").Value = "  -5.72%  "
$ws.Range("E12").Value = "  -5.75%  "
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").Value = "3.015.62"
$ws.Range("E14").Value = "  -4.23%  "
$ws.Range("D15").Value = "58.352.01"
$ws.Range("E15").Value = "  -4.93%  "
$ws.Range("D16").Value = "21.06"
$ws.Range("E16").Value = "  -5.44%  "
$ws.Range("E17").Value = "  -5.17%  "
$ws.Range("D18").Value = "2.577.21"
$ws.Range("E18").Value = "  -4.42%  "
$ws.Range("D19").Value = "4.54"
$ws.Range("E19").Value = "  -5.68%  "
$ws.Range("D20").Value = "342.81"
$ws.Range("E20").Value = "  -4.18%  "
$ws.Range("D21").Value = "10.27"
$ws.Range("E21").Value = "  -4.86%  "
$ws.Range("D22").Value = "6.05"
$ws.Range("E22").Value = "  -5.32%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "60.52"
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("D25").Value = "0.418"
$ws.Range("E25").Value = "  -4.29%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "2.677.67"
$ws.Range("E27").Value = "  -4.11%  "
$ws.Range("E28").Value = "  -6.52%  "
$ws.Range("D29").Value = "0.0₃0812"
$ws.Range("E29").Value = "  -7.05%  "
$ws.Range("D30").Value = "7.00"
$ws.Range("E30").Value = "  -6.13%  "
$ws.Range("D32").Value = "6.09"
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("D33").Value = "18.77"
$ws.Range("E33").Value = "  -4.58%  "
$ws.Range("D34").Value = "149.32"
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("D35").Value = "1.54"
$ws.Range("E35").Value = "  -6.16%  "
$ws.Range("D36").Value = "0.946"
$ws.Range("E36").Value = "  +6.08%  "
$ws.Range("D37").Value = "3.97"
$ws.Range("E37").Value = "  -5.25%  "
$ws.Range("E38").Value = "  -7.41%  "
$ws.Range("E39").Value = "  -7.67%  "
$ws.Range("D40").Value = "36.10"
$ws.Range("D41").Value = "291.52"
$ws.Range("E41").Value = "  -5.40%  "
$ws.Range("E42").Value = "  -7.59%  "
$ws.Range("E43").Value = "  -6.45%  "
$ws.Range("D44").Value = "0.0993"
$ws.Range("E44").Value = "  -3.40%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "0.608"
$ws.Range("E46").Value = "  -7.77%  "
$ws.Range("E47").Value = "  -6.05%  "
$ws.Range("D48").Value = "19.21"
$ws.Range("E48").Value = "  -8.84%  "
$ws.Range("D49").Value = "10.27"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0228"
$ws.Range("E50").Value = "  -6.05%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "4.67"
$ws.Range("E51").Value = "  -7.71%  "
